$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-08 Monday" "2024-04-09 Tuesday"

Replace-Text "389÷3=129, 2" "437÷3=145, 2"
Replace-Text "779÷9=86, 5" "794÷2=397, 0"
Replace-Text "788÷6=131, 2" "384÷4=96, 0"
Replace-Text "843÷3=281, 0" "669÷3=223, 0"
Replace-Text "397÷7=56, 5" "827÷6=137, 5"

Replace-Text "926÷4=231, 2" "423÷2=211, 1"
Replace-Text "946÷3=315, 1" "854÷8=106, 6"
Replace-Text "820÷4=205, 0" "476÷5=95, 1"
Replace-Text "301÷4=75, 1" "580÷2=290, 0"
Replace-Text "127÷4=31, 3" "789÷2=394, 1"

Replace-Text "989÷7=141, 2" "462÷9=51, 3"
Replace-Text "922÷7=131, 5" "518÷4=129, 2"
Replace-Text "269÷3=89, 2" "451÷9=50, 1"
Replace-Text "648÷7=92, 4" "160÷2=80, 0"
Replace-Text "542÷2=271, 0" "197÷2=98, 1"

Replace-Text "577÷3=192, 1" "890÷3=296, 2"
Replace-Text "261÷9=29, 0" "188÷2=94, 0"
Replace-Text "962÷2=481, 0" "434÷9=48, 2"
Replace-Text "783÷2=391, 1" "450÷9=50, 0"
Replace-Text "895÷7=127, 6" "815÷7=116, 3"

Replace-Text "997÷8=124, 5" "841÷4=210, 1"
Replace-Text "724÷3=241, 1" "499÷4=124, 3"
Replace-Text "476÷8=59, 4" "395÷5=79, 0"
Replace-Text "517÷7=73, 6" "553÷4=138, 1"
Replace-Text "245÷5=49, 0" "680÷3=226, 2"
